$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '243.95'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '25.07'
$ws.Range('D3').ClearFormats()
$ws.Range('B4').Value = 'LEO'
$ws.Range('C4').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '3.501'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '3LEOLEO'
$ws.Range('B5').Value = 'HuobiToken'
$ws.Range('C5').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '5.186'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '4HuobiTokenHT'
$ws.Range('B6').Value = 'Cronos'
$ws.Range('C6').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.05749'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '5CronosCRO'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.504'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.115'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '7GateTokenGT'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8086'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '8MXTokenMX'
$ws.Range('B10').Value = 'FTXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8404'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '9FTXTokenFTT'
$ws.Range('B11').Value = 'One'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.009662'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '10OneONE'
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.1339'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '11WazirXWRX'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.06959'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '12MandalaExchangeTokenMDX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.02827'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.09368'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.001511'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.006247'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '16TigerCashTCH'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3174'
$ws.Range('D19').ClearFormats()
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1301'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.739'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04659'
$ws.Range('D23').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.001237'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.004271'
$ws.Range('D26').ClearFormats()
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0001985'
$ws.Range('D28').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03610'
$ws.Range('D40').ClearFormats()
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1049'
$ws.Range('D42').ClearFormats()
$ws.Range('E43').Value = '42CEJICEJIBestin24h'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.007331'
$ws.Range('D44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005299'
$ws.Range('D45').ClearFormats()
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.2829'
$ws.Range('D47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.002276'
$ws.Range('D48').ClearFormats()
